$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new data row was inserted at row 44 ("Fruta / hortaliza, semanal" update),
# pushing the former rows 44..172 down to 45..173.
$ws.Rows.Item(44).Insert()

# Populate the newly inserted row 44 with its data.
$ws.Cells.Item(44, 1).Value2  = 5
$ws.Cells.Item(44, 2).Value2  = "Macroferia Regional de Talca"
$ws.Cells.Item(44, 3).Value2  = "Maule"
$ws.Cells.Item(44, 4).Value2  = 44949
$ws.Cells.Item(44, 5).Value2  = 7
$ws.Cells.Item(44, 6).Value2  = 100112030
$ws.Cells.Item(44, 7).Value2  = "Poroto granado"
$ws.Cells.Item(44, 8).Value2  = "Sin especificar"
$ws.Cells.Item(44, 9).Value2  = "Primera"
$ws.Cells.Item(44, 10).Value2 = 150
$ws.Cells.Item(44, 11).Value2 = 45000
$ws.Cells.Item(44, 12).Value2 = 45000
$ws.Cells.Item(44, 13).Value2 = 45000
$ws.Cells.Item(44, 14).Value2 = "`$/saco 25 kilos"
$ws.Cells.Item(44, 15).Value2 = "Región del Maule"
$ws.Cells.Item(44, 16).Value2 = 1800
$ws.Cells.Item(44, 17).Value2 = 25
$ws.Cells.Item(44, 18).Value2 = "Hortaliza"
